# Add a new inventory row for the 10G Ethernet card (row 23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "10G Ethernet card"
$ws.Range("B23").Value = "10 GIGABIT ETHERNET PCIE ADAPTER CARD FOR DESKTOP - ETTUS RESEARCH"
$ws.Range("C23").Value = "HC2494813"
$ws.Range("E23").Value = "San Diego - Palomar"

# Carry over the formatting used by the rest of the table (style index "3")
# for the Item and Network Location columns, matching the existing rows.
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("E22").Copy()
$ws.Range("E23").PasteSpecial(-4122) # xlPasteFormats

# Update the active selection to just below the newly added row
[void]$ws.Range("E24").Select()
